$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted at row 35 ("Haba" / Feria
# Lagunitas de Puerto Montt), pushing the existing rows 35-120 down to
# 36-121 (dimension grows from A1:R120 to A1:R121).
$ws.Rows("35:35").Insert()

$ws.Range("A35").Value = 4
$ws.Range("B35").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C35").Value = 'Los Lagos'
$ws.Range("D35").Value = 44883
$ws.Range("E35").Value = 10
$ws.Range("F35").Value = 100112026
$ws.Range("G35").Value = 'Haba'
$ws.Range("H35").Value = 'Sin especificar'
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = 15000
$ws.Range("N35").Value = '$/saco 25 kilos'
$ws.Range("O35").Value = 'Región Metropolitana'
$ws.Range("P35").Value = 600
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = 'Hortaliza'
